$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("admin inv 2", "admininv2@gmail.com", "12346", "pkr", "Civil", "professor"),
    @("admin inv 3", "admininv3@gmail.com", "12347", "ktm", "Computer", "teacher"),
    @("admin inv 4", "admininv4@gmail.com", "12348", "pkr", "Civil", "professor"),
    @("admin inv 5", "admininv5@gmail.com", "12349", "ktm", "Computer", "teacher"),
    @("admin inv 6", "admininv6@gmail.com", "12350", "pkr", "Civil", "professor"),
    @("admin inv 7", "admininv7@gmail.com", "12351", "ktm", "Computer", "teacher"),
    @("admin inv 8", "admininv8@gmail.com", "12352", "pkr", "Civil", "professor"),
    @("admin inv 9", "admininv9@gmail.com", "12353", "ktm", "Computer", "teacher")
)

$ws.Range("A1").Copy()

$rowIndex = 2
foreach ($rec in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $rec[0]
    $ws.Cells.Item($rowIndex, 2).Value = $rec[1]
    $ws.Cells.Item($rowIndex, 3).Value = "'" + $rec[2]
    $ws.Cells.Item($rowIndex, 4).Value = $rec[3]
    $ws.Cells.Item($rowIndex, 5).Value = $rec[4]
    $ws.Cells.Item($rowIndex, 6).Value = $rec[5]

    $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122)

    $rowIndex = $rowIndex + 1
}
